$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calibrations")
$ws.Range("A1").Value = "test"
